# Updated self grading sheet.
$wb = $excel.ActiveWorkbook

# "Basic Game rubric" sheet: update scores
$basic = $wb.Worksheets.Item("Basic Game rubric")
$basic.Range("B2").Value = 3   # Camera: 2 -> 3
$basic.Range("B4").Value = 3   # Interactions: 2 -> 3
$basic.Range("B5").Value = 2   # Game implementation: 3 -> 2
$basic.Range("B6").Value = 1   # HUD and UI: 0 -> 1

# "Game extras" sheet: note that "Load level from a file" used an external library
$extras = $wb.Worksheets.Item("Game extras")
$extras.Range("C4").Value = "External library"
